$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "261k01"
$ws.Range("A3").Value = "261k02"

$ws.Range("A4").Select()
